# Update column G ("K" - strikeouts) values for rows 2-24 of Sheet1.
# These values were regenerated upstream (switching the source stat from
# "Strike#" to "K"); write the new literal values into the worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 4
    3  = 3
    4  = 3
    5  = 7
    6  = 6
    7  = 1
    8  = 3
    9  = 4
    10 = 5
    11 = 8
    12 = 6
    13 = 11
    14 = 0
    15 = 6
    16 = 4
    17 = 8
    18 = 5
    19 = 6
    20 = 6
    21 = 7
    22 = 7
    23 = 6
    24 = 3
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
